# Auto-generated edit script applying numeric corrections to Gungnir_Profits workbook
# Updates cached market-price / profit figures on several crafting-class sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 442.6
$ws.Range("I33").Value = 358.0909
$ws.Range("J33").Value = 675
$ws.Range("K33").Value = 358.0909
$ws.Range("L33").Value = 675
$ws.Range("M33").Value = -129.0909
$ws.Range("N33").Value = -1133
# Row 47
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21944
# Row 74
$ws.Range("H74").Value = 6930.9
$ws.Range("I74").Value = 3746.5
$ws.Range("K74").Value = 3746.5
$ws.Range("M74").Value = -2810.5
# Row 77
$ws.Range("H77").Value = 6930.9
$ws.Range("I77").Value = 3746.5
$ws.Range("K77").Value = 18732.5
$ws.Range("M77").Value = -14052.5
# Row 86
$ws.Range("H86").Value = 26475
$ws.Range("J86").Value = 2500
$ws.Range("L86").Value = 2500
$ws.Range("N86").Value = -4746
# Row 89
$ws.Range("H89").Value = 26475
$ws.Range("J89").Value = 2500
$ws.Range("L89").Value = 12500
$ws.Range("N89").Value = -23732
# Row 134
$ws.Range("H134").Value = 47857.145
$ws.Range("J134").Value = 47857.145
$ws.Range("L134").Value = 47857.145
$ws.Range("N134").Value = -57997.145
# Row 138
$ws.Range("H138").Value = 3263.5557
$ws.Range("I138").Value = 1351.3438
$ws.Range("J138").Value = 4176.8506
$ws.Range("K138").Value = 4054.0314
$ws.Range("L138").Value = 12530.5518
$ws.Range("M138").Value = 1085.9686
$ws.Range("N138").Value = -22810.5518

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1723.3462
$ws.Range("I61").Value = 1599.8518
$ws.Range("J61").Value = 1856.72
$ws.Range("K61").Value = 1599.8518
$ws.Range("L61").Value = 1856.72
$ws.Range("M61").Value = -1387.8518
$ws.Range("N61").Value = -2280.72
# Row 136
$ws.Range("H136").Value = 1723.3462
$ws.Range("I136").Value = 1599.8518
$ws.Range("J136").Value = 1856.72
$ws.Range("K136").Value = 4799.555399999999
$ws.Range("L136").Value = 5570.16
$ws.Range("M136").Value = -2249.555399999999
$ws.Range("N136").Value = -10670.16

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2586952.8
$ws.Range("I86").Value = 4472.4
$ws.Range("J86").Value = 5815053
$ws.Range("K86").Value = 4472.4
$ws.Range("L86").Value = 5815053
$ws.Range("M86").Value = -3349.4
$ws.Range("N86").Value = -5817299
# Row 89
$ws.Range("H89").Value = 2586952.8
$ws.Range("I89").Value = 4472.4
$ws.Range("J89").Value = 5815053
$ws.Range("K89").Value = 22362
$ws.Range("L89").Value = 29075265
$ws.Range("M89").Value = -16746
$ws.Range("N89").Value = -29086497

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 5500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 8000
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = -388
$ws.Range("N4").Value = -8224
# Row 59
$ws.Range("H59").Value = 47780
$ws.Range("J59").Value = 47780
$ws.Range("L59").Value = 47780
$ws.Range("N59").Value = -50070
# Row 108
$ws.Range("H108").Value = 55833.332
$ws.Range("J108").Value = 55833.332
$ws.Range("L108").Value = 55833.332
$ws.Range("N108").Value = -63513.332
# Row 141
$ws.Range("H141").Value = 57600
$ws.Range("I141").Value = 17000
$ws.Range("J141").Value = 64981.816
$ws.Range("K141").Value = 17000
$ws.Range("L141").Value = 64981.816
$ws.Range("M141").Value = -11820
$ws.Range("N141").Value = -75341.81599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 101.666664
$ws.Range("I4").Value = 62
$ws.Range("K4").Value = 186
$ws.Range("M4").Value = -74
# Row 121
$ws.Range("H121").Value = 8334634.5
$ws.Range("I121").Value = 700
$ws.Range("J121").Value = 11112613
$ws.Range("K121").Value = 2100
$ws.Range("L121").Value = 33337839
$ws.Range("M121").Value = -790
$ws.Range("N121").Value = -33340459
# Row 122
$ws.Range("H122").Value = 19235376
$ws.Range("I122").Value = 41667064
$ws.Range("J122").Value = 8214.214
$ws.Range("K122").Value = 375003576
$ws.Range("L122").Value = 73927.92600000001
$ws.Range("M122").Value = -375001126
$ws.Range("N122").Value = -78827.92600000001
# Row 134
$ws.Range("H134").Value = 1854.0625
$ws.Range("I134").Value = 1644.4667
$ws.Range("J134").Value = 4998
$ws.Range("K134").Value = 4933.4001
$ws.Range("L134").Value = 14994
$ws.Range("M134").Value = 136.5999000000002
$ws.Range("N134").Value = -25134
# Row 136
$ws.Range("H136").Value = 7579281
$ws.Range("I136").Value = 25001304
$ws.Range("J136").Value = 4487.913
$ws.Range("K136").Value = 75003912
$ws.Range("L136").Value = 13463.739
$ws.Range("M136").Value = -74998812
$ws.Range("N136").Value = -23663.739

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 12000
$ws.Range("J5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("N5").Value = -12224
# Row 39
$ws.Range("H39").Value = 28261
$ws.Range("J39").Value = 28261
$ws.Range("L39").Value = 28261
$ws.Range("N39").Value = -29325
# Row 80
$ws.Range("H80").Value = 6670314
$ws.Range("I80").Value = 4567.1113
$ws.Range("J80").Value = 16668934
$ws.Range("K80").Value = 4567.1113
$ws.Range("L80").Value = 16668934
$ws.Range("M80").Value = -3569.1113
$ws.Range("N80").Value = -16670930
# Row 83
$ws.Range("H83").Value = 6670314
$ws.Range("I83").Value = 4567.1113
$ws.Range("J83").Value = 16668934
$ws.Range("K83").Value = 22835.5565
$ws.Range("L83").Value = 83344670
$ws.Range("M83").Value = -17843.5565
$ws.Range("N83").Value = -83354654
# Row 112
$ws.Range("H112").Value = 37183.855
$ws.Range("J112").Value = 37183.855
$ws.Range("L112").Value = 37183.855
$ws.Range("N112").Value = -39399.855
# Row 126
$ws.Range("H126").Value = 4522.9414
$ws.Range("I126").Value = 5020
$ws.Range("J126").Value = 4315.8335
$ws.Range("K126").Value = 15060
$ws.Range("L126").Value = 12947.5005
$ws.Range("M126").Value = -12590
$ws.Range("N126").Value = -17887.5005
# Row 132
$ws.Range("H132").Value = 4361.4326
$ws.Range("I132").Value = 1186.5
$ws.Range("J132").Value = 10222.846
$ws.Range("K132").Value = 3559.5
$ws.Range("L132").Value = 30668.538
$ws.Range("M132").Value = -1029.5
$ws.Range("N132").Value = -35728.538

$ws = $wb.Worksheets.Item("LTW")
# Row 104
$ws.Range("H104").Value = 41275
$ws.Range("J104").Value = 41275
$ws.Range("L104").Value = 41275
$ws.Range("N104").Value = -48263
# Row 111
$ws.Range("H111").Value = 40450
$ws.Range("J111").Value = 40450
$ws.Range("L111").Value = 40450
$ws.Range("N111").Value = -48630
# Row 132
$ws.Range("H132").Value = 13517543
$ws.Range("I132").Value = 20409570
$ws.Range("J132").Value = 9168.32
$ws.Range("K132").Value = 61228710
$ws.Range("L132").Value = 27504.96
$ws.Range("M132").Value = -61226180
$ws.Range("N132").Value = -32564.96

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 12361.111
$ws.Range("I2").Value = 38125
$ws.Range("K2").Value = 38125
$ws.Range("M2").Value = -38013
# Row 41
$ws.Range("H41").Value = 4229.4287
$ws.Range("I41").Value = 2500
$ws.Range("J41").Value = 4517.6665
$ws.Range("K41").Value = 2500
$ws.Range("L41").Value = 4517.6665
$ws.Range("M41").Value = -2110
$ws.Range("N41").Value = -5297.6665
# Row 45
$ws.Range("H45").Value = 11166.667
$ws.Range("I45").Value = 3000
$ws.Range("J45").Value = 12800
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 12800
$ws.Range("M45").Value = -2509
$ws.Range("N45").Value = -13782
# Row 108
$ws.Range("H108").Value = 36648
$ws.Range("J108").Value = 36648
$ws.Range("L108").Value = 36648
$ws.Range("N108").Value = -44328
